$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "iso_image_file" column (L) entirely - both its header and its
# one data value - shifting the sheet's used range back to column K.
$ws.Range("L1:L2").EntireColumn.Delete() | Out-Null

# Update the virtualization software string for the "Kali" row (row 3).
$ws.Range("D3").Value = "VirtualBox Version 7.1.2 r164945 (Qt6.5.3)"

# Fix the OS instance string for the "Kali" row (row 3) to match row 2.
$ws.Range("E3").Value = "Kali 2024.2"

# Move the selection to D3 (and scroll the view back to column A).
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D3").Select() | Out-Null
